$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 10 (8 rows). The projection data series moves to
# start 8 quarters later (dates/values that used to live in rows 11+ now
# occupy rows 3+), and every subsequent row shifts up by 8, shrinking the
# used range from B2:D118 down to B2:D110.
$ws.Rows("3:10").Delete()

# Rename the worksheet from "Sheet1" to "data".
$ws.Name = "data"
